# "Update outputs for 2024"
# Populate the first four benchmark rows (Day / C# average runtime in
# seconds) that the RuntimesChart bar chart reads from ($A$3:$A$20 /
# $B$3:$B$20), and leave the selection where Excel would land after typing
# the values in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0.00241846
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0.00499994
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 0.00228852
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 0.00462774

$ws.Range("A3:B6").Select() | Out-Null
